# Update cryptos list Price (D) / Volume(1h) (E) columns per the latest data pull.
# Note: several "Price" values are plain decimals (e.g. "205.68") that Excel would
# otherwise auto-coerce to numbers; a leading apostrophe forces text entry (matching
# the source data's text type), and resetting .Style back to "Normal" afterwards
# avoids leaving a stray quote-prefix format on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.895.50'
$ws.Range("E2").Value = '  +0.17%  '

$ws.Range("D3").Value = '1.547.08'
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").Value = '''205.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.05%  '

$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").Value = '''21.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.21%  '

$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("D12").Value = '1.767.42'
$ws.Range("E12").Value = '  -0.95%  '

$ws.Range("D13").Value = '1.554.42'
$ws.Range("E13").Value = '  -0.58%  '

$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("D16").Value = '26.871.00'
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").Value = '''61.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.49%  '

$ws.Range("D18").Value = '''213.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("D19").Value = '0.0₃0681'
$ws.Range("E19").Value = '  +0.16%  '

$ws.Range("D20").Value = '''7.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.09%  '

$ws.Range("E21").Value = '  +0.29%  '

$ws.Range("E22").Value = '  -2.57%  '

$ws.Range("E23").Value = '  -0.29%  '

$ws.Range("E24").Value = '  -3.28%  '

$ws.Range("D25").Value = '''152.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("E26").Value = '  -0.84%  '

$ws.Range("E27").Value = '  -0.61%  '

$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("E31").Value = '  -1.01%  '

$ws.Range("E32").Value = '  +1.85%  '

$ws.Range("D33").Value = '1.362.54'
$ws.Range("E33").Value = '  -2.76%  '

$ws.Range("E34").Value = '  +0.77%  '

$ws.Range("E35").Value = '  +0.38%  '

$ws.Range("E36").Value = '  +5.37%  '

$ws.Range("E37").Value = '  +0.35%  '

$ws.Range("D38").Value = '''0.0165'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").Value = '''0.519'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.56%  '

$ws.Range("D40").Value = '''0.805'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.67%  '

$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("D42").Value = '''5.59'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.70%  '

$ws.Range("D43").Value = '''0.986'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.92%  '

$ws.Range("E44").Value = '  +1.53%  '

$ws.Range("D45").Value = '''63.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.18%  '

$ws.Range("E46").Value = '  -2.14%  '

$ws.Range("D47").Value = '1.681.66'
$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("D48").Value = '''86.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.35%  '

$ws.Range("D49").Value = '''0.0506'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.49%  '

$ws.Range("D50").Value = '0.0₇0964'
$ws.Range("E50").Value = '  -1.19%  '

$ws.Range("D51").Value = '''0.0948'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.07%  '
